$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FirstSheet")

$ws.Range("A10").Value = "Abhi_0"
$ws.Range("B10").Value = "Abhi_1"
$ws.Range("C10").Value = "Abhi_2"

$ws.Range("A11").Value = "Abhi_0"
$ws.Range("B11").Value = "Abhi_1"
$ws.Range("C11").Value = "Abhi_2"
